$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new EUR->ARS quote row (row 44) to the table.
$row = 44

# Column A holds a literal date string like "2025-09-27". Force the cell to
# Text format BEFORE assigning the value so Excel doesn't auto-convert the
# ISO-looking string into a date serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-27"

$ws.Cells.Item($row, 2).Value = "21:21:16"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,623.5666"
